# Clean up the leading spaces on the state names in Sheet2 (column A),
# then leave Sheet2 selected/active as the last user action.

$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item(2)

# Trim the one-off "Wisconsin" cell first (this is what the author's manual
# edit apparently touched individually), then clean up the rest of the
# column in a single batch operation. The order here matters because it
# determines the order in which new shared-string entries get created.
$wisconsinCell = $ws2.Range("A21")
$wisconsinCell.Value = $wisconsinCell.Value2.Trim()

$rng = $ws2.Range("A2:A51")
$vals = $rng.Value2
for ($i = 1; $i -le 50; $i++) {
    $orig = $vals[$i, 1]
    if ($orig -ne $null) {
        $vals[$i, 1] = $orig.Trim()
    }
}
$rng.Value2 = $vals

# Make Sheet2 the active sheet and select the full state-name column,
# matching the final selection/active-tab state left behind by the edit.
$ws2.Activate()
$ws2.Range("A2:A51").Select()
